$d = $word.ActiveDocument

# --- Insertion 1 ---
# A new, empty "ListParagraph" paragraph (inheriting the black/text1 run-mark
# color already used by the paragraph above it) is added right before the
# "Model Selection" heading.
$needle1 = "so next, we use word embedding to help improving model performance."
$replace1 = $needle1 + "^p"
$d.Content.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2) | Out-Null

# --- Insertion 2 ---
# Right after the paragraph ending "...to occur." two new paragraphs are
# added:
#   1. "First of all,  we use the LSTM model with one-hot method"
#   2. an empty paragraph
# both styled as "ListParagraph" with the black/text1 run color, followed by
# the pre-existing empty paragraph.
$needle2 = "It predicts a word given in the user input and then each of the next words is predicted using the probability of likelihood of that word to occur."
$newText = "First of all, " + " we use the " + "LSTM model with one-hot method"
$replace2 = $needle2 + "^p" + $newText + "^p"
$d.Content.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2) | Out-Null
